$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1442.0847
$ws.Range("I15").Value = 1442.0847
$ws.Range("K15").Value = 4326.2541
$ws.Range("M15").Value = -4157.2541

$ws.Range("H41").Value = 563.8333
$ws.Range("I41").Value = 358
$ws.Range("J41").Value = 975.5
$ws.Range("K41").Value = 358
$ws.Range("L41").Value = 975.5
$ws.Range("M41").Value = 82
$ws.Range("N41").Value = -1855.5

$ws.Range("H43").Value = 6299.933
$ws.Range("I43").Value = 4153.615
$ws.Range("K43").Value = 4153.615
$ws.Range("M43").Value = -4084.615

$ws.Range("H62").Value = 136061.75
$ws.Range("I62").Value = 136061.75
$ws.Range("K62").Value = 136061.75
$ws.Range("M62").Value = -135437.75

$ws.Range("H65").Value = 136061.75
$ws.Range("I65").Value = 136061.75
$ws.Range("K65").Value = 680308.75
$ws.Range("M65").Value = -677188.75

$ws.Range("H74").Value = 42800.8
$ws.Range("I74").Value = 42800.8
$ws.Range("K74").Value = 42800.8
$ws.Range("M74").Value = -41864.8

$ws.Range("H77").Value = 42800.8
$ws.Range("I77").Value = 42800.8
$ws.Range("K77").Value = 214004
$ws.Range("M77").Value = -209324

$ws.Range("H101").Value = 1624.8334
$ws.Range("I101").Value = 1124.75
$ws.Range("K101").Value = 3374.25
$ws.Range("M101").Value = -1752.25

$ws.Range("H106").Value = 41489.96
$ws.Range("I106").Value = 45619.74
$ws.Range("K106").Value = 45619.74
$ws.Range("M106").Value = -44988.74

$ws.Range("H107").Value = 1027.8334
$ws.Range("I107").Value = 336.85715
$ws.Range("J107").Value = 1995.2
$ws.Range("K107").Value = 336.85715
$ws.Range("L107").Value = 1995.2
$ws.Range("M107").Value = 1583.14285
$ws.Range("N107").Value = -5835.2

$ws.Range("H112").Value = 4921.846
$ws.Range("J112").Value = 4990.5
$ws.Range("L112").Value = 14971.5
$ws.Range("N112").Value = -17187.5

$ws.Range("H121").Value = 1498.4615
$ws.Range("J121").Value = 1498.4615
$ws.Range("L121").Value = 4495.3845
$ws.Range("N121").Value = -7989.3845

$ws.Range("H132").Value = 50009836
$ws.Range("J132").Value = 14832
$ws.Range("L132").Value = 44496
$ws.Range("N132").Value = -49556

$ws.Range("H135").Value = 6120.148
$ws.Range("J135").Value = 13866.333
$ws.Range("L135").Value = 124796.997
$ws.Range("N135").Value = -129866.997

$ws.Range("H137").Value = 1491.3462
$ws.Range("I137").Value = 1161
$ws.Range("J137").Value = 1518.875
$ws.Range("K137").Value = 3483
$ws.Range("L137").Value = 4556.625
$ws.Range("M137").Value = -933
$ws.Range("N137").Value = -9656.625

$ws.Range("H138").Value = 268503.75
$ws.Range("I138").Value = 2033.0625
$ws.Range("K138").Value = 6099.1875
$ws.Range("M138").Value = -959.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 215.13043
$ws.Range("I5").Value = 204.85
$ws.Range("J5").Value = 283.66666
$ws.Range("K5").Value = 204.85
$ws.Range("L5").Value = 283.66666
$ws.Range("M5").Value = -92.85
$ws.Range("N5").Value = -507.66666

$ws.Range("H6").Value = 10500000
$ws.Range("J6").Value = 1000000
$ws.Range("L6").Value = 1000000
$ws.Range("N6").Value = -1000346

$ws.Range("H32").Value = 2360.7031
$ws.Range("I32").Value = 2415.5518
$ws.Range("K32").Value = 2415.5518
$ws.Range("M32").Value = -2128.5518

$ws.Range("H74").Value = 2914.28
$ws.Range("I74").Value = 3690.8
$ws.Range("J74").Value = 2396.6
$ws.Range("K74").Value = 3690.8
$ws.Range("L74").Value = 2396.6
$ws.Range("M74").Value = -2816.8
$ws.Range("N74").Value = -4144.6

$ws.Range("H77").Value = 2914.28
$ws.Range("I77").Value = 3690.8
$ws.Range("J77").Value = 2396.6
$ws.Range("K77").Value = 18454
$ws.Range("L77").Value = 11983
$ws.Range("M77").Value = -14086
$ws.Range("N77").Value = -20719

$ws.Range("H97").Value = 3008.8965
$ws.Range("I97").Value = 1628.56
$ws.Range("J97").Value = 11636
$ws.Range("K97").Value = 1628.56
$ws.Range("L97").Value = 11636
$ws.Range("M97").Value = -1132.56
$ws.Range("N97").Value = -12628

$ws.Range("H132").Value = 2150.3677
$ws.Range("I132").Value = 2136.4575
$ws.Range("J132").Value = 2241.5557
$ws.Range("K132").Value = 6409.372499999999
$ws.Range("L132").Value = 6724.6671
$ws.Range("M132").Value = -3879.372499999999
$ws.Range("N132").Value = -11784.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 215.13043
$ws.Range("I4").Value = 204.85
$ws.Range("J4").Value = 283.66666
$ws.Range("K4").Value = 204.85
$ws.Range("L4").Value = 283.66666
$ws.Range("M4").Value = -89.85
$ws.Range("N4").Value = -513.66666

$ws.Range("H107").Value = 3096.8
$ws.Range("I107").Value = 3124.923
$ws.Range("K107").Value = 3124.923
$ws.Range("M107").Value = -1204.923

$ws.Range("H123").Value = 53500
$ws.Range("I123").Value = 47000
$ws.Range("J123").Value = 60000
$ws.Range("K123").Value = 47000
$ws.Range("L123").Value = 60000
$ws.Range("M123").Value = -42100
$ws.Range("N123").Value = -69800

$ws.Range("H134").Value = 2772.2917
$ws.Range("I134").Value = 2772.2917
$ws.Range("K134").Value = 8316.875100000001
$ws.Range("M134").Value = -5781.875100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1250.25
$ws.Range("I16").Value = 1250.25
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1250.25
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -963.25
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value = 3699.889
$ws.Range("I31").Value = 2542
$ws.Range("J31").Value = 4094.0637
$ws.Range("K31").Value = 2542
$ws.Range("L31").Value = 4094.0637
$ws.Range("M31").Value = -2247
$ws.Range("N31").Value = -4684.063700000001

$ws.Range("H34").Value = 3699.889
$ws.Range("I34").Value = 2542
$ws.Range("J34").Value = 4094.0637
$ws.Range("K34").Value = 2542
$ws.Range("L34").Value = 4094.0637
$ws.Range("M34").Value = -2340
$ws.Range("N34").Value = -4498.063700000001

$ws.Range("H58").Value = 3619.7778
$ws.Range("I58").Value = 5694.5
$ws.Range("K58").Value = 5694.5
$ws.Range("M58").Value = -5491.5

$ws.Range("H93").Value = 12696.4
$ws.Range("J93").Value = 49969
$ws.Range("L93").Value = 49969
$ws.Range("N93").Value = -53713

$ws.Range("H99").Value = 5118.1934
$ws.Range("I99").Value = 4945.3213
$ws.Range("K99").Value = 4945.3213
$ws.Range("M99").Value = -3447.3213

$ws.Range("H105").Value = 8445.857
$ws.Range("I105").Value = 8557.333
$ws.Range("K105").Value = 8557.333
$ws.Range("M105").Value = -6810.333000000001

$ws.Range("H113").Value = 1250.25
$ws.Range("I113").Value = 1250.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1250.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 919.75
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 5118.1934
$ws.Range("I126").Value = 4945.3213
$ws.Range("K126").Value = 14835.9639
$ws.Range("M126").Value = -12365.9639

$ws.Range("H132").Value = 2106.2144
$ws.Range("I132").Value = 1986.5416
$ws.Range("J132").Value = 2824.25
$ws.Range("K132").Value = 5959.6248
$ws.Range("L132").Value = 8472.75
$ws.Range("M132").Value = -3429.6248
$ws.Range("N132").Value = -13532.75

$ws.Range("H134").Value = 6325.4053
$ws.Range("I134").Value = 6003.5625
$ws.Range("K134").Value = 18010.6875
$ws.Range("M134").Value = -15475.6875

$ws.Range("H136").Value = 3619.7778
$ws.Range("I136").Value = 5694.5
$ws.Range("K136").Value = 17083.5
$ws.Range("M136").Value = -14533.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9989.75
$ws.Range("I3").Value = 3320.8333
$ws.Range("J3").Value = 29996.5
$ws.Range("K3").Value = 9962.499899999999
$ws.Range("L3").Value = 89989.5
$ws.Range("M3").Value = -9850.499899999999
$ws.Range("N3").Value = -90213.5

$ws.Range("H5").Value = 521.6226
$ws.Range("I5").Value = 476.375
$ws.Range("J5").Value = 529.6667
$ws.Range("K5").Value = 1429.125
$ws.Range("L5").Value = 1589.0001
$ws.Range("M5").Value = -1317.125
$ws.Range("N5").Value = -1813.0001

$ws.Range("H7").Value = 259.7857
$ws.Range("I7").Value = 209.77777
$ws.Range("J7").Value = 349.8
$ws.Range("K7").Value = 629.33331
$ws.Range("L7").Value = 1049.4
$ws.Range("M7").Value = -517.33331
$ws.Range("N7").Value = -1273.4

$ws.Range("H34").Value = 1308.12
$ws.Range("J34").Value = 400
$ws.Range("L34").Value = 1200
$ws.Range("N34").Value = -1368

$ws.Range("H39").Value = 1794.5714
$ws.Range("I39").Value = 1510.3334
$ws.Range("K39").Value = 4531.0002
$ws.Range("M39").Value = -4237.0002

$ws.Range("H68").Value = 2517.0908
$ws.Range("J68").Value = 2720.5
$ws.Range("L68").Value = 8161.5
$ws.Range("N68").Value = -9783.5

$ws.Range("H71").Value = 2517.0908
$ws.Range("J71").Value = 2720.5
$ws.Range("L71").Value = 24484.5
$ws.Range("N71").Value = -32596.5

$ws.Range("H113").Value = 490.86487
$ws.Range("J113").Value = 570.38464
$ws.Range("L113").Value = 1711.15392
$ws.Range("N113").Value = -6051.15392

$ws.Range("H121").Value = 2328.4644
$ws.Range("I121").Value = 297.30768
$ws.Range("J121").Value = 4088.8
$ws.Range("K121").Value = 891.92304
$ws.Range("L121").Value = 12266.4
$ws.Range("M121").Value = 418.07696
$ws.Range("N121").Value = -14886.4

$ws.Range("H129").Value = 3673.5557
$ws.Range("J129").Value = 5173.75
$ws.Range("L129").Value = 15521.25
$ws.Range("N129").Value = -25521.25

$ws.Range("H131").Value = 4790492.5
$ws.Range("I131").Value = 20793.5
$ws.Range("K131").Value = 62380.5
$ws.Range("M131").Value = -57340.5

$ws.Range("H132").Value = 1989.4445
$ws.Range("I132").Value = 1950
$ws.Range("K132").Value = 17550
$ws.Range("M132").Value = -15020

$ws.Range("H135").Value = 521.6226
$ws.Range("I135").Value = 476.375
$ws.Range("J135").Value = 529.6667
$ws.Range("K135").Value = 4287.375
$ws.Range("L135").Value = 4767.0003
$ws.Range("M135").Value = -1752.375
$ws.Range("N135").Value = -9837.0003

$ws.Range("H137").Value = 12556247
$ws.Range("J137").Value = 83960.46
$ws.Range("L137").Value = 251881.38
$ws.Range("N137").Value = -262081.38

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 54301.96
$ws.Range("J43").Value = 90692
$ws.Range("L43").Value = 90692
$ws.Range("N43").Value = -90994

$ws.Range("H57").Value = 64027.5
$ws.Range("J57").Value = 84036.664
$ws.Range("L57").Value = 84036.664
$ws.Range("N57").Value = -85676.664

$ws.Range("H59").Value = 12300
$ws.Range("J59").Value = 12300
$ws.Range("L59").Value = 12300
$ws.Range("N59").Value = -13466

$ws.Range("H70").Value = 7316.5
$ws.Range("I70").Value = 5216.067
$ws.Range("K70").Value = 5216.067
$ws.Range("M70").Value = -4946.067

$ws.Range("H73").Value = 7316.5
$ws.Range("I73").Value = 5216.067
$ws.Range("K73").Value = 5216.067
$ws.Range("M73").Value = -4280.067

$ws.Range("H80").Value = 5093.357
$ws.Range("I80").Value = 3030.9
$ws.Range("K80").Value = 3030.9
$ws.Range("M80").Value = -2032.9

$ws.Range("H83").Value = 5093.357
$ws.Range("I83").Value = 3030.9
$ws.Range("K83").Value = 15154.5
$ws.Range("M83").Value = -10162.5

$ws.Range("H97").Value = 1309.6666
$ws.Range("I97").Value = 824.8571
$ws.Range("K97").Value = 824.8571
$ws.Range("M97").Value = -328.8570999999999

$ws.Range("H132").Value = 6852.625
$ws.Range("I132").Value = 7759.75
$ws.Range("K132").Value = 23279.25
$ws.Range("M132").Value = -20749.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4776
$ws.Range("I16").Value = 1450.5
$ws.Range("K16").Value = 1450.5
$ws.Range("M16").Value = -1280.5

$ws.Range("H46").Value = 9969.385
$ws.Range("I46").Value = 4027.75
$ws.Range("K46").Value = 4027.75
$ws.Range("M46").Value = -3839.75

$ws.Range("H55").Value = 980.1875
$ws.Range("I55").Value = 291.35294
$ws.Range("J55").Value = 1760.8667
$ws.Range("K55").Value = 291.35294
$ws.Range("L55").Value = 1760.8667
$ws.Range("M55").Value = -118.35294
$ws.Range("N55").Value = -2106.8667

$ws.Range("H93").Value = 1035.6818
$ws.Range("I93").Value = 897.17645
$ws.Range("J93").Value = 1506.6
$ws.Range("K93").Value = 897.17645
$ws.Range("L93").Value = 1506.6
$ws.Range("M93").Value = 350.82355
$ws.Range("N93").Value = -4002.6

$ws.Range("H100").Value = 4100.3
$ws.Range("J100").Value = 4500
$ws.Range("L100").Value = 4500
$ws.Range("N100").Value = -5582

$ws.Range("H127").Value = 79666.5
$ws.Range("J127").Value = 79666.5
$ws.Range("L127").Value = 79666.5
$ws.Range("N127").Value = -89586.5

$ws.Range("H132").Value = 4436.1665
$ws.Range("I132").Value = 4472.1514
$ws.Range("J132").Value = 4304.222
$ws.Range("K132").Value = 13416.4542
$ws.Range("L132").Value = 12912.666
$ws.Range("M132").Value = -10886.4542
$ws.Range("N132").Value = -17972.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 83446.9
$ws.Range("I62").Value = 160395.2
$ws.Range("K62").Value = 160395.2
$ws.Range("M62").Value = -159771.2

$ws.Range("H65").Value = 83446.9
$ws.Range("I65").Value = 160395.2
$ws.Range("K65").Value = 801976
$ws.Range("M65").Value = -798856

$ws.Range("H80").Value = 40000
$ws.Range("J80").Value = 40000
$ws.Range("L80").Value = 40000
$ws.Range("N80").Value = -41996

$ws.Range("H83").Value = 40000
$ws.Range("J83").Value = 40000
$ws.Range("L83").Value = 120000
$ws.Range("N83").Value = -129984

$ws.Range("H100").Value = 2357.75
$ws.Range("J100").Value = 3521.111
$ws.Range("L100").Value = 7042.222
$ws.Range("N100").Value = -8124.222

$ws.Range("H103").Value = 35950
$ws.Range("J103").Value = 35950
$ws.Range("L103").Value = 35950
$ws.Range("N103").Value = -38294

$ws.Range("H113").Value = 6945671
$ws.Range("I113").Value = 16667687
$ws.Range("J113").Value = 1373.8572
$ws.Range("K113").Value = 50003061
$ws.Range("L113").Value = 4121.571599999999
$ws.Range("M113").Value = -50000891
$ws.Range("N113").Value = -8461.5716

$ws.Range("H132").Value = 2176.347
$ws.Range("I132").Value = 1620.6875
$ws.Range("K132").Value = 4862.0625
$ws.Range("M132").Value = -2332.0625
